$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 text with new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 12820.13 pesos`n✅ 12820.13 pesos = 3.35 = 966.04 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 298
$wsTasas.Range("O10").Value = 3820.4
$wsTasas.Range("N12").Value = 3822
$wsTasas.Range("O12").Value = 288
